$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# --- Footer (Footers.Item(1)) : Pearson logo, docPr id="2" -> rename image1.png to image2.png ---
$ftr1 = $sec.Footers.Item(1)
$ishp1 = $ftr1.Range.InlineShapes.Item(1)
$shape1 = $ishp1.ConvertToShape()
$shape1.Name = "image2.png"
$shape1.ConvertToInlineShape()

# --- Footer (Footers.Item(2)) : Pearson logo, docPr id="3" -> rename image1.png to image2.png ---
$ftr2 = $sec.Footers.Item(2)
$ishp2 = $ftr2.Range.InlineShapes.Item(1)
$shape2 = $ishp2.ConvertToShape()
$shape2.Name = "image2.png"
$shape2.ConvertToInlineShape()

# --- Header (Headers.Item(2)) : BTec logo, docPr id="1" -> rename image2.jpg to image1.jpg ---
$hdr2 = $sec.Headers.Item(2)
$ishp3 = $hdr2.Range.InlineShapes.Item(1)
$shape3 = $ishp3.ConvertToShape()
$shape3.Name = "image1.jpg"
$shape3.ConvertToInlineShape()

$d.Save()
